$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "necessary"
$ws.Range("B16").Value = "77777K000"

$ws.Range("B16").Font.Color = 0
$ws.Range("B16").Font.Name = "Arial"
$ws.Range("B16").Font.Size = 10

$ws.Range("A16").Select()
